# Fix column A (File_name / "pdf file name" description rows).
# Rows 3-10 currently hold the short placeholder text "pdf file name";
# they should hold the same descriptive text already used in row 2:
# "pdf file name (text)."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$correctText = $ws.Range("A2").Value2

for ($r = 3; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $correctText
}

# Move the active selection to B8 (matches the saved view state).
[void]$ws.Range("B8").Select()
